# Refresh the cryptos price/volume snapshot (GitHub Actions bot update).
# Price-column (D) values are written with a leading apostrophe so Excel's
# smart-entry parser stores them as text instead of silently converting
# number-looking strings (e.g. "1.00", "7.70") into numeric values and
# dropping the trailing zero; the Style reset that follows clears the
# quote-prefix formatting so the cell style matches the original (unstyled)
# text cells. Volume-column (E) values already contain non-numeric
# characters ("%", padding spaces) so they round-trip as text unaided.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'50.963.81"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = "'2.945.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'379.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = "'101.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('D10').Value = "'36.13"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = "'0.0851"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = "'3.403.98"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').Value = "'18.26"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').Value = "'7.70"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.01%  '
$ws.Range('D16').Value = "'12.13"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +68.15%  '
$ws.Range('D17').Value = "'2.945.77"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = "'0.995"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = "'50.914.67"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -4.88%  '
$ws.Range('D21').Value = "'12.40"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = "'69.45"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('D24').Value = "'266.80"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = "'3.24"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +13.15%  '
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = "'7.08"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.60%  '
$ws.Range('D29').Value = "'25.61"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  -3.25%  '
$ws.Range('D31').Value = "'0.108"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('D32').Value = "'10.04"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('D33').Value = "'50.54"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = "'33.51"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = "'3.12"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.12%  '
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').Value = "'16.56"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('D42').Value = "'2.50"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('D43').Value = "'120.43"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('D44').Value = "'21.36"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('D45').Value = "'3.45"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.22%  '
$ws.Range('E46').Value = '  -2.11%  '
$ws.Range('D47').Value = "'2.33"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = "'2.011.36"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  -4.68%  '
$ws.Range('E50').Value = '  -6.21%  '
$ws.Range('D51').Value = "'5.30"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.91%  '
